$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "61.932.50"
Set-TextValue $ws "E2" "  -0.70%  "
Set-TextValue $ws "D3" "2.442.35"
Set-TextValue $ws "E3" "  +0.58%  "
Set-TextValue $ws "E4" "  -0.15%  "
Set-TextValue $ws "D5" "579.32"
Set-TextValue $ws "E5" "  +1.12%  "
Set-TextValue $ws "D6" "141.12"
Set-TextValue $ws "E6" "  -1.16%  "
Set-TextValue $ws "E7" "  +0.03%  "
Set-TextValue $ws "D8" "0.529"
Set-TextValue $ws "E8" "  +0.38%  "
Set-TextValue $ws "D9" "2.437.25"
Set-TextValue $ws "E9" "  +0.65%  "
Set-TextValue $ws "D10" "0.109"
Set-TextValue $ws "E10" "  +2.45%  "
Set-TextValue $ws "E11" "  +2.66%  "
Set-TextValue $ws "D12" "5.17"
Set-TextValue $ws "E12" "  +0.18%  "
Set-TextValue $ws "D13" "0.339"
Set-TextValue $ws "E13" "  -2.24%  "
Set-TextValue $ws "D14" "25.91"
Set-TextValue $ws "E14" "  -1.10%  "
Set-TextValue $ws "E15" "  +0.17%  "
Set-TextValue $ws "D16" "2.893.41"
Set-TextValue $ws "E16" "  +0.27%  "
Set-TextValue $ws "D17" "61.775.39"
Set-TextValue $ws "E17" "  -0.69%  "
Set-TextValue $ws "D18" "2.442.83"
Set-TextValue $ws "E18" "  +0.24%  "
Set-TextValue $ws "D19" "10.59"
Set-TextValue $ws "E19" "  -3.85%  "
Set-TextValue $ws "D20" "7.20"
Set-TextValue $ws "E20" "  +1.89%  "
Set-TextValue $ws "D21" "325.12"
Set-TextValue $ws "E21" "  -0.59%  "
Set-TextValue $ws "D22" "4.05"
Set-TextValue $ws "E22" "  -1.27%  "
Set-TextValue $ws "E23" "  +0.14%  "
Set-TextValue $ws "E24" "  -5.05%  "
Set-TextValue $ws "D25" "65.15"
Set-TextValue $ws "E25" "  -0.49%  "
Set-TextValue $ws "D26" "9.16"
Set-TextValue $ws "E26" "  +1.89%  "
Set-TextValue $ws "D27" "590.94"
Set-TextValue $ws "E27" "  -5.61%  "
Set-TextValue $ws "E28" "  -0.09%  "
Set-TextValue $ws "E29" "  +0.20%  "
Set-TextValue $ws "D30" "0.0₃0939"
Set-TextValue $ws "E30" "  -0.84%  "
Set-TextValue $ws "E31" "  -1.49%  "
Set-TextValue $ws "E32" "  -2.48%  "
Set-TextValue $ws "D33" "1.87"
Set-TextValue $ws "E33" "  +0.32%  "
Set-TextValue $ws "E34" "  -2.32%  "
Set-TextValue $ws "E35" "  +0.03%  "
Set-TextValue $ws "D36" "4.77"
Set-TextValue $ws "E36" "  -2.85%  "
Set-TextValue $ws "B37" "PolygonEcosystemToken"
Set-TextValue $ws "C37" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws "D37" "0.373"
Set-TextValue $ws "E37" "  -0.07%  "
Set-TextValue $ws "B38" "Monero"
Set-TextValue $ws "C38" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D38" "152.23"
Set-TextValue $ws "E38" "  +3.74%  "
Set-TextValue $ws "E39" "  -2.50%  "
Set-TextValue $ws "D40" "18.31"
Set-TextValue $ws "E40" "  -0.16%  "
Set-TextValue $ws "D41" "5.18"
Set-TextValue $ws "E41" "  -0.01%  "
Set-TextValue $ws "D42" "42.96"
Set-TextValue $ws "E42" "  +1.49%  "
Set-TextValue $ws "D43" "1.00"
Set-TextValue $ws "E43" "  +0.08%  "
Set-TextValue $ws "E44" "  -2.43%  "
Set-TextValue $ws "D45" "2.37"
Set-TextValue $ws "E45" "  -2.20%  "
Set-TextValue $ws "B46" "BabyDogeCoin"
Set-TextValue $ws "C46" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D46" "0.0₆0271"
Set-TextValue $ws "E46" "  +20.12%  "
Set-TextValue $ws "B47" "Aave"
Set-TextValue $ws "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D47" "140.73"
Set-TextValue $ws "E47" "  -2.07%  "
Set-TextValue $ws "D48" "3.58"
Set-TextValue $ws "E48" "  -2.34%  "
Set-TextValue $ws "E49" "  +0.42%  "
Set-TextValue $ws "D50" "0.0511"
Set-TextValue $ws "E50" "  -1.36%  "
Set-TextValue $ws "D51" "19.53"
Set-TextValue $ws "E51" "  +0.83%  "

Write-Host "Done applying changes"
